$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two new weekly records are inserted at rows 343-344 which pushes the
# existing rows 343-387 down to 345-389 (dimension grows from R387 to R389).
$ws.Rows("343:344").Insert()

# New row 343
$ws.Cells.Item(343, 1).Value = 7
$ws.Cells.Item(343, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(343, 3).Value = "Ñuble"
$ws.Cells.Item(343, 4).Value = 44984
$ws.Cells.Item(343, 5).Value = 16
$ws.Cells.Item(343, 6).Value = 100114013
$ws.Cells.Item(343, 7).Value = "Zanahoria"
$ws.Cells.Item(343, 8).Value = "Sin especificar"
$ws.Cells.Item(343, 9).Value = "Primera"
$ws.Cells.Item(343, 10).Value = 240
$ws.Cells.Item(343, 11).Value = 6500
$ws.Cells.Item(343, 12).Value = 7000
$ws.Cells.Item(343, 13).Value = 6750
$ws.Cells.Item(343, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(343, 15).Value = "Región de Ñuble"
$ws.Cells.Item(343, 16).Value = 338
$ws.Cells.Item(343, 17).Value = 20
$ws.Cells.Item(343, 18).Value = "Hortaliza"

# New row 344
$ws.Cells.Item(344, 1).Value = 7
$ws.Cells.Item(344, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(344, 3).Value = "Ñuble"
$ws.Cells.Item(344, 4).Value = 44984
$ws.Cells.Item(344, 5).Value = 16
$ws.Cells.Item(344, 6).Value = 100114013
$ws.Cells.Item(344, 7).Value = "Zanahoria"
$ws.Cells.Item(344, 8).Value = "Sin especificar"
$ws.Cells.Item(344, 9).Value = "Segunda"
$ws.Cells.Item(344, 10).Value = 100
$ws.Cells.Item(344, 11).Value = 5500
$ws.Cells.Item(344, 12).Value = 5500
$ws.Cells.Item(344, 13).Value = 5500
$ws.Cells.Item(344, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(344, 15).Value = "Región de Ñuble"
$ws.Cells.Item(344, 16).Value = 275
$ws.Cells.Item(344, 17).Value = 20
$ws.Cells.Item(344, 18).Value = "Hortaliza"
